$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(76, 8).Value = 6194.3657
$ws.Cells.Item(76, 9).Value = 7721.0835
$ws.Cells.Item(76, 10).Value = 4039
$ws.Cells.Item(76, 11).Value = 7721.0835
$ws.Cells.Item(76, 12).Value = 4039
$ws.Cells.Item(76, 13).Value = -7406.0835
$ws.Cells.Item(76, 14).Value = -4669
$ws.Cells.Item(79, 8).Value = 6194.3657
$ws.Cells.Item(79, 9).Value = 7721.0835
$ws.Cells.Item(79, 10).Value = 4039
$ws.Cells.Item(79, 11).Value = 7721.0835
$ws.Cells.Item(79, 12).Value = 4039
$ws.Cells.Item(79, 13).Value = -6629.0835
$ws.Cells.Item(79, 14).Value = -6223
$ws.Cells.Item(86, 8).Value = 2162.25
$ws.Cells.Item(86, 9).Value = 2291.4546
$ws.Cells.Item(86, 10).Value = 2052.923
$ws.Cells.Item(86, 11).Value = 2291.4546
$ws.Cells.Item(86, 12).Value = 2052.923
$ws.Cells.Item(86, 13).Value = -1168.4546
$ws.Cells.Item(86, 14).Value = -4298.923
$ws.Cells.Item(89, 8).Value = 2162.25
$ws.Cells.Item(89, 9).Value = 2291.4546
$ws.Cells.Item(89, 10).Value = 2052.923
$ws.Cells.Item(89, 11).Value = 11457.273
$ws.Cells.Item(89, 12).Value = 10264.615
$ws.Cells.Item(89, 13).Value = -5841.273000000001
$ws.Cells.Item(89, 14).Value = -21496.615
$ws.Cells.Item(106, 8).Value = 63494068
$ws.Cells.Item(106, 9).Value = 41668376
$ws.Cells.Item(106, 10).Value = 76925260
$ws.Cells.Item(106, 11).Value = 41668376
$ws.Cells.Item(106, 12).Value = 76925260
$ws.Cells.Item(106, 13).Value = -41667745
$ws.Cells.Item(106, 14).Value = -76926522
$ws.Cells.Item(132, 8).Value = 1689.1428
$ws.Cells.Item(132, 9).Value = 1985.6
$ws.Cells.Item(132, 10).Value = 948
$ws.Cells.Item(132, 11).Value = 5956.799999999999
$ws.Cells.Item(132, 12).Value = 2844
$ws.Cells.Item(132, 13).Value = -3426.799999999999
$ws.Cells.Item(132, 14).Value = -7904
$ws.Cells.Item(137, 8).Value = 1259.7742
$ws.Cells.Item(137, 9).Value = 1235.8214
$ws.Cells.Item(137, 10).Value = 1483.3334
$ws.Cells.Item(137, 11).Value = 3707.4642
$ws.Cells.Item(137, 12).Value = 4450.0002
$ws.Cells.Item(137, 13).Value = -1157.4642
$ws.Cells.Item(137, 14).Value = -9550.0002

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(132, 8).Value = 3073.6047
$ws.Cells.Item(132, 9).Value = 1580.5
$ws.Cells.Item(132, 10).Value = 4959.6313
$ws.Cells.Item(132, 11).Value = 4741.5
$ws.Cells.Item(132, 12).Value = 14878.8939
$ws.Cells.Item(132, 13).Value = -2211.5
$ws.Cells.Item(132, 14).Value = -19938.8939

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(99, 8).Value = 166668140
$ws.Cells.Item(99, 9).Value = 250001000
$ws.Cells.Item(99, 10).Value = 2435
$ws.Cells.Item(99, 11).Value = 250001000
$ws.Cells.Item(99, 12).Value = 2435
$ws.Cells.Item(99, 13).Value = -249999502
$ws.Cells.Item(99, 14).Value = -5431

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 30118.324
$ws.Cells.Item(31, 9).Value = 47441.547
$ws.Cells.Item(31, 10).Value = 4710.933
$ws.Cells.Item(31, 11).Value = 47441.547
$ws.Cells.Item(31, 12).Value = 4710.933
$ws.Cells.Item(31, 13).Value = -47146.547
$ws.Cells.Item(31, 14).Value = -5300.933
$ws.Cells.Item(34, 8).Value = 30118.324
$ws.Cells.Item(34, 9).Value = 47441.547
$ws.Cells.Item(34, 10).Value = 4710.933
$ws.Cells.Item(34, 11).Value = 47441.547
$ws.Cells.Item(34, 12).Value = 4710.933
$ws.Cells.Item(34, 13).Value = -47239.547
$ws.Cells.Item(34, 14).Value = -5114.933
$ws.Cells.Item(58, 8).Value = 1457.1786
$ws.Cells.Item(58, 9).Value = 1066.7894
$ws.Cells.Item(58, 10).Value = 2281.3333
$ws.Cells.Item(58, 11).Value = 1066.7894
$ws.Cells.Item(58, 12).Value = 2281.3333
$ws.Cells.Item(58, 13).Value = -863.7893999999999
$ws.Cells.Item(58, 14).Value = -2687.3333
$ws.Cells.Item(62, 8).Value = 5453.636
$ws.Cells.Item(62, 9).Value = 6600
$ws.Cells.Item(62, 10).Value = 4078
$ws.Cells.Item(62, 11).Value = 6600
$ws.Cells.Item(62, 12).Value = 4078
$ws.Cells.Item(62, 13).Value = -5976
$ws.Cells.Item(62, 14).Value = -5326
$ws.Cells.Item(65, 8).Value = 5453.636
$ws.Cells.Item(65, 9).Value = 6600
$ws.Cells.Item(65, 10).Value = 4078
$ws.Cells.Item(65, 11).Value = 33000
$ws.Cells.Item(65, 12).Value = 20390
$ws.Cells.Item(65, 13).Value = -29880
$ws.Cells.Item(65, 14).Value = -26630
$ws.Cells.Item(136, 8).Value = 1457.1786
$ws.Cells.Item(136, 9).Value = 1066.7894
$ws.Cells.Item(136, 10).Value = 2281.3333
$ws.Cells.Item(136, 11).Value = 3200.3682
$ws.Cells.Item(136, 12).Value = 6843.999899999999
$ws.Cells.Item(136, 13).Value = -650.3681999999999
$ws.Cells.Item(136, 14).Value = -11943.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 354076.72
$ws.Cells.Item(5, 9).Value = 756
$ws.Cells.Item(5, 10).Value = 751562.5
$ws.Cells.Item(5, 11).Value = 2268
$ws.Cells.Item(5, 12).Value = 2254687.5
$ws.Cells.Item(5, 13).Value = -2156
$ws.Cells.Item(5, 14).Value = -2254911.5
$ws.Cells.Item(124, 8).Value = 4999.7144
$ws.Cells.Item(124, 9).Value = 1500
$ws.Cells.Item(124, 10).Value = 5583
$ws.Cells.Item(124, 11).Value = 4500
$ws.Cells.Item(124, 12).Value = 16749
$ws.Cells.Item(124, 13).Value = 410
$ws.Cells.Item(124, 14).Value = -26569
$ws.Cells.Item(131, 8).Value = 2000970
$ws.Cells.Item(131, 9).Value = 10000460
$ws.Cells.Item(131, 10).Value = 1097.6
$ws.Cells.Item(131, 11).Value = 30001380
$ws.Cells.Item(131, 12).Value = 3292.8
$ws.Cells.Item(131, 13).Value = -29996340
$ws.Cells.Item(131, 14).Value = -13372.8
$ws.Cells.Item(135, 8).Value = 354076.72
$ws.Cells.Item(135, 9).Value = 756
$ws.Cells.Item(135, 10).Value = 751562.5
$ws.Cells.Item(135, 11).Value = 6804
$ws.Cells.Item(135, 12).Value = 6764062.5
$ws.Cells.Item(135, 13).Value = -4269
$ws.Cells.Item(135, 14).Value = -6769132.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 2593.4138
$ws.Cells.Item(80, 9).Value = 2497.647
$ws.Cells.Item(80, 10).Value = 2729.0833
$ws.Cells.Item(80, 11).Value = 2497.647
$ws.Cells.Item(80, 12).Value = 2729.0833
$ws.Cells.Item(80, 13).Value = -1499.647
$ws.Cells.Item(80, 14).Value = -4725.0833
$ws.Cells.Item(83, 8).Value = 2593.4138
$ws.Cells.Item(83, 9).Value = 2497.647
$ws.Cells.Item(83, 10).Value = 2729.0833
$ws.Cells.Item(83, 11).Value = 12488.235
$ws.Cells.Item(83, 12).Value = 13645.4165
$ws.Cells.Item(83, 13).Value = -7496.235000000001
$ws.Cells.Item(83, 14).Value = -23629.4165
$ws.Cells.Item(132, 8).Value = 3209
$ws.Cells.Item(132, 9).Value = 3987.75
$ws.Cells.Item(132, 10).Value = 2819.625
$ws.Cells.Item(132, 11).Value = 11963.25
$ws.Cells.Item(132, 12).Value = 8458.875
$ws.Cells.Item(132, 13).Value = -9433.25
$ws.Cells.Item(132, 14).Value = -13518.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 47620028
$ws.Cells.Item(46, 9).Value = 83333976
$ws.Cells.Item(46, 10).Value = 1433.3334
$ws.Cells.Item(46, 11).Value = 83333976
$ws.Cells.Item(46, 12).Value = 1433.3334
$ws.Cells.Item(46, 13).Value = -83333788
$ws.Cells.Item(46, 14).Value = -1809.3334
$ws.Cells.Item(132, 8).Value = 11757065
$ws.Cells.Item(132, 9).Value = 19103498
$ws.Cells.Item(132, 10).Value = 2773.8
$ws.Cells.Item(132, 11).Value = 57310494
$ws.Cells.Item(132, 12).Value = 8321.400000000001
$ws.Cells.Item(132, 13).Value = -57307964
$ws.Cells.Item(132, 14).Value = -13381.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 4733.6665
$ws.Cells.Item(62, 9).Value = 4467.3335
$ws.Cells.Item(62, 10).Value = 5000
$ws.Cells.Item(62, 11).Value = 4467.3335
$ws.Cells.Item(62, 12).Value = 5000
$ws.Cells.Item(62, 13).Value = -3843.3335
$ws.Cells.Item(62, 14).Value = -6248
$ws.Cells.Item(65, 8).Value = 4733.6665
$ws.Cells.Item(65, 9).Value = 4467.3335
$ws.Cells.Item(65, 10).Value = 5000
$ws.Cells.Item(65, 11).Value = 22336.6675
$ws.Cells.Item(65, 12).Value = 25000
$ws.Cells.Item(65, 13).Value = -19216.6675
$ws.Cells.Item(65, 14).Value = -31240
$ws.Cells.Item(132, 8).Value = 1367.0303
$ws.Cells.Item(132, 9).Value = 1075.4
$ws.Cells.Item(132, 10).Value = 2278.375
$ws.Cells.Item(132, 11).Value = 3226.2
$ws.Cells.Item(132, 12).Value = 6835.125
$ws.Cells.Item(132, 13).Value = -696.2000000000003
$ws.Cells.Item(132, 14).Value = -11895.125
